$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "Metal" -> "Metal (mark only)"
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Metal (mark only)"

# ---------------------------------------------------------------------------
# 2. New row 11: Scrapbook Paper (machine block 2 only)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Scrapbook Paper"
$ws.Range("H11").Value = 26.1
$ws.Range("I11").Value = 15
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = "off"

# ---------------------------------------------------------------------------
# 3. New row 12: Quilted Glass (machine block 3 only)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Quilted Glass"
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 19
$ws.Range("O12").Value = 165

# ---------------------------------------------------------------------------
# 4. Header row (row 4) + data row (row 5): add a 4th "Z-Axis" column to
#    each of the three machine blocks. Blocks 2 and 3 shift one column to
#    the right. Propagate the block's cell format to the new column BEFORE
#    moving values, so the style carries over like Excel's own "insert
#    column"/fill behaviour would.
# ---------------------------------------------------------------------------
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)

$ws.Range("K4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("L4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("M4").Copy()
$ws.Range("O4").PasteSpecial(-4122)

$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Copy()
$ws.Range("J5").PasteSpecial(-4122)

# Now move the actual values (the original cells still hold their original
# content - PasteSpecial above only touched formats).
$r4_power2 = $ws.Range("G4").Value()
$r4_speed2 = $ws.Range("H4").Value()
$r4_ppi2   = $ws.Range("I4").Value()
$r4_power3 = $ws.Range("K4").Value()
$r4_speed3 = $ws.Range("L4").Value()
$r4_ppi3   = $ws.Range("M4").Value()

$r5_power2 = $ws.Range("G5").Value()
$r5_speed2 = $ws.Range("H5").Value()
$r5_ppi2   = $ws.Range("I5").Value()
$r5_power3 = $ws.Range("K5").Value()
$r5_speed3 = $ws.Range("L5").Value()

$ws.Range("H4").Value = $r4_power2
$ws.Range("I4").Value = $r4_speed2
$ws.Range("J4").Value = $r4_ppi2
$ws.Range("M4").Value = $r4_power3
$ws.Range("N4").Value = $r4_speed3
$ws.Range("O4").Value = $r4_ppi3

$ws.Range("H5").Value = $r5_power2
$ws.Range("I5").Value = $r5_speed2
$ws.Range("J5").Value = $r5_ppi2
$ws.Range("N5").Value = $r5_power3
$ws.Range("O5").Value = $r5_speed3

# New "Z-Axis" header cells, matching font/colour of their block.
$ws.Range("F4").Value = "Z-Axis"
$ws.Range("F4").Font.Bold = $true

$ws.Range("K4").Value = "Z-Axis"
$ws.Range("K4").Font.Bold = $true
$ws.Range("K4").Font.Color = 255

$ws.Range("P4").Value = "Z-Axis"
$ws.Range("P4").Font.Bold = $true
$ws.Range("P4").Font.Color = 12611584

# Drop the now-unused old columns.
$ws.Range("G4").Clear()
$ws.Range("L4").Clear()
$ws.Range("G5").Clear()
$ws.Range("K5").Clear()
$ws.Range("L5").Clear()

# ---------------------------------------------------------------------------
# 5. Row 8 (1/8 Cipboard row): shift block 2 data right (no styling on the
#    original cells, so a plain value move is enough). Values are written
#    as literals (rather than round-tripped through Value()) to avoid
#    reformatting 3.8 as a long binary-float string.
# ---------------------------------------------------------------------------
$ws.Range("H8").Value = 95
$ws.Range("I8").Value = 3.8
$ws.Range("J8").Value = 500

$ws.Range("G8").Clear()

# ---------------------------------------------------------------------------
# 6. Finish row 12 (Quilted Glass): Z-Axis value "on"
# ---------------------------------------------------------------------------
$ws.Range("P12").Value = "on"

# ---------------------------------------------------------------------------
# 7. New row 13: Hardboard (all three machine blocks, all Z-Axis "on")
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Hardboard"
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = 3.75
$ws.Range("E13").Value = 600
$ws.Range("F13").Value = "on"
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 12
$ws.Range("J13").Value = 600
$ws.Range("K13").Value = "on"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 3.75
$ws.Range("O13").Value = 600
$ws.Range("P13").Value = "on"

# ---------------------------------------------------------------------------
# 8. New row 14: Birch Ply (section header, like "Cipboard")
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Birch Ply"

# ---------------------------------------------------------------------------
# 9. New rows 15-17: thickness labels (1/16, 1/8, 1/4) styled like rows 7-9
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "1/16"
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").HorizontalAlignment = -4152

$ws.Range("A16").Value = "1/8"
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").HorizontalAlignment = -4152

$ws.Range("A17").Value = "1/4"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").HorizontalAlignment = -4152

$ws.Range("C17").Value = 100
$ws.Range("D17").Value = 70
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = "on"
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 2.5
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = "on"

# ---------------------------------------------------------------------------
# 10. New row 18: Glass (Soda Ash)
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "Glass (Soda Ash)"
$ws.Range("C18").Value = 100
$ws.Range("D18").Value = 53
$ws.Range("E18").Value = 500
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = 19
$ws.Range("O18").Value = 165

# ---------------------------------------------------------------------------
# 11. View state: scroll back to show column A, select A19
# ---------------------------------------------------------------------------
$ws.Range("A19").Select()
